$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B93: it was stored as text ("4"); make it a real number, as in the rest
# of the politeness_score column.
$ws.Cells.Item(93, 2).Value = 4

# Add new row 94 with the additional annotation.
$ws.Cells.Item(94, 1).Value = "Ruilin"

# B94 mirrors the original (buggy) B93 entry: a numeric-looking value stored
# as text. Force text storage, then restore the default style so no stray
# cell-level style reference is left behind.
$ws.Cells.Item(94, 2).NumberFormat = "@"
$ws.Cells.Item(94, 2).Value = "3"
$ws.Cells.Item(94, 2).Style = "Normal"

$ws.Cells.Item(94, 3).Value = "无"
$ws.Cells.Item(94, 4).Value = "DFT"
$ws.Cells.Item(94, 5).Value = "MET"
$ws.Cells.Item(94, 6).Value = "8fd9d1eb-d55a-4b83-a989-0f77ecdd42b7"
$ws.Cells.Item(94, 7).Value = "r1CE9GWR-_annotated.xlsx"
$ws.Cells.Item(94, 8).Value = "Moreover, the discussion of supervised and unsupervised paradigms is utterly unconvincing, especially in light of the above comment on minimum-distance estimation underlying both of these paradigms."
